$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B66 currently holds text "2" as an inline string; convert to a real number.
$ws.Cells.Item(66, 2).Value = 2

# Add new row 67 with the additional annotation data.
$ws.Cells.Item(67, 1).Value = "Ying Tang"

# B67 is "4" stored as TEXT (not a number) in the target, so force the
# cell to text format before writing the value, then restore the default
# "Normal" style so no visible formatting change remains.
$ws.Cells.Item(67, 2).NumberFormat = "@"
$ws.Cells.Item(67, 2).Value = "4"
$ws.Cells.Item(67, 2).Style = "Normal"

$ws.Cells.Item(67, 3).Value = " It would be nice"
$ws.Cells.Item(67, 4).Value = "SUG"
$ws.Cells.Item(67, 5).Value = "RES"
$ws.Cells.Item(67, 6).Value = "b01bb119-e44b-4008-9381-38115d7c20f9"
$ws.Cells.Item(67, 7).Value = "mugzy2nI-Ayi1_annotated.xlsx"
$ws.Cells.Item(67, 8).Value = "It would be nice to have more explanation of the significance of beating SignalP."
